# cboc_signin_sheet.xlsx — widen the Month/Year header row on both sheets
# from its old width out to column AE, and bump the month from FEBRUARY to
# APRIL 2021.
#
# Both worksheets have row 1 as a single merged "Month/Year: ..." banner
# cell. We grow the merge by inserting new columns just before the current
# last column of the merged range (EntireColumn.Insert on a merged range
# automatically extends the merge area, and the newly inserted cells pick
# up the formatting of the column immediately to their left, while the
# pushed-right original last cell keeps its own original formatting). This
# preserves the existing cell styles (header / middle / end borders)
# instead of Excel's Merge()/MergeCells=$true path, which always
# recomputes formatting for the whole merged block from the top-left cell.

$wb = $excel.ActiveWorkbook

$newMonthYear = "Month/Year: APRIL 2021"

# --- Sheet "1-15": merged range currently A1:O1, needs to become A1:AE1 ---
# O (15) -> AE (31): 16 new columns, inserted immediately before O.
$ws1 = $wb.Worksheets.Item("1-15")
$ws1.Range("O1:AD1").EntireColumn.Insert()
$ws1.Range("A1").Value = $newMonthYear

# --- Sheet "16-End": merged range currently A1:M1, needs to become A1:AE1 ---
# M (13) -> AE (31): 18 new columns, inserted immediately before M.
$ws2 = $wb.Worksheets.Item("16-End")
$ws2.Range("M1:AD1").EntireColumn.Insert()
$ws2.Range("A1").Value = $newMonthYear
